$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 7500
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 7500
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 7500
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -8468
$ws.Range("H70").Value = 1812.125
$ws.Range("J70").Value = 1812.125
$ws.Range("L70").Value = 5436.375
$ws.Range("N70").Value = -5976.375
$ws.Range("H73").Value = 1812.125
$ws.Range("J73").Value = 1812.125
$ws.Range("L73").Value = 5436.375
$ws.Range("N73").Value = -7308.375
$ws.Range("H86").Value = 2667
$ws.Range("I86").Value = 2501.5
$ws.Range("J86").Value = 2998
$ws.Range("K86").Value = 2501.5
$ws.Range("L86").Value = 2998
$ws.Range("M86").Value = -1378.5
$ws.Range("N86").Value = -5244
$ws.Range("H89").Value = 2667
$ws.Range("I89").Value = 2501.5
$ws.Range("J89").Value = 2998
$ws.Range("K89").Value = 12507.5
$ws.Range("L89").Value = 14990
$ws.Range("M89").Value = -6891.5
$ws.Range("N89").Value = -26222
$ws.Range("H100").Value = 1050
$ws.Range("J100").Value = 1200
$ws.Range("L100").Value = 1200
$ws.Range("N100").Value = -2282
$ws.Range("H106").Value = 1495
$ws.Range("I106").Value = 1495
$ws.Range("K106").Value = 1495
$ws.Range("M106").Value = -864

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2089.3
$ws.Range("I2").Value = 2099.2222
$ws.Range("K2").Value = 2099.2222
$ws.Range("M2").Value = -1986.2222
$ws.Range("H10").Value = 4
$ws.Range("I10").Value = 4
$ws.Range("K10").Value = 4
$ws.Range("M10").Value = 166
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H32").Value = 834
$ws.Range("I32").Value = 834
$ws.Range("K32").Value = 834
$ws.Range("M32").Value = -547
$ws.Range("H63").Value = 7725
$ws.Range("I63").Value = 7725
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 7725
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -7039
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 7725
$ws.Range("I66").Value = 7725
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 38625
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -35193
$ws.Range("N66").ClearContents()
$ws.Range("H97").Value = 998
$ws.Range("I97").Value = 998
$ws.Range("K97").Value = 998
$ws.Range("M97").Value = -502
$ws.Range("H116").Value = 2089.3
$ws.Range("I116").Value = 2099.2222
$ws.Range("K116").Value = 2099.2222
$ws.Range("M116").Value = 194.7777999999998

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2089.3
$ws.Range("I3").Value = 2099.2222
$ws.Range("K3").Value = 2099.2222
$ws.Range("M3").Value = -1985.2222
$ws.Range("H36").Value = 772.6667
$ws.Range("I36").Value = 772.6667
$ws.Range("K36").Value = 772.6667
$ws.Range("M36").Value = -238.6667
$ws.Range("H94").Value = 3333
$ws.Range("I94").Value = 2749.5
$ws.Range("J94").Value = 4500
$ws.Range("K94").Value = 2749.5
$ws.Range("L94").Value = 4500
$ws.Range("M94").Value = -2298.5
$ws.Range("N94").Value = -5402

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 245.38889
$ws.Range("I2").Value = 274.375
$ws.Range("K2").Value = 1646.25
$ws.Range("M2").Value = -1533.25
$ws.Range("H38").Value = 88
$ws.Range("J38").Value = 90
$ws.Range("L38").Value = 270
$ws.Range("N38").Value = -964
$ws.Range("H109").Value = 2700
$ws.Range("I109").Value = 2700
$ws.Range("K109").Value = 8100
$ws.Range("M109").Value = -7060
$ws.Range("H131").Value = 1933.4872
$ws.Range("I131").Value = 5666.3335
$ws.Range("J131").Value = 1622.4166
$ws.Range("K131").Value = 16999.0005
$ws.Range("L131").Value = 4867.2498
$ws.Range("M131").Value = -11959.0005
$ws.Range("N131").Value = -14947.2498

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H31").Value = 800
$ws.Range("I31").Value = 1100
$ws.Range("J31").Value = 200
$ws.Range("K31").Value = 1100
$ws.Range("L31").Value = 200
$ws.Range("M31").Value = -808
$ws.Range("N31").Value = -784
$ws.Range("H37").Value = 800
$ws.Range("I37").Value = 1100
$ws.Range("J37").Value = 200
$ws.Range("K37").Value = 1100
$ws.Range("L37").Value = 200
$ws.Range("M37").Value = -823
$ws.Range("N37").Value = -754
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 280
$ws.Range("I20").Value = 280
$ws.Range("K20").Value = 280
$ws.Range("M20").Value = -54
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H93").Value = 2494
$ws.Range("I93").Value = 2494
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 2494
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -1246
$ws.Range("N93").ClearContents()
$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32774
$ws.Range("H132").Value = 11166.667
$ws.Range("I132").Value = 9500
$ws.Range("K132").Value = 28500
$ws.Range("M132").Value = -25970

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 4
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("H69").Value = 27499.75
$ws.Range("J69").Value = 27499.75
$ws.Range("L69").Value = 27499.75
$ws.Range("N69").Value = -28997.75
$ws.Range("H72").Value = 27499.75
$ws.Range("J72").Value = 27499.75
$ws.Range("L72").Value = 82499.25
$ws.Range("N72").Value = -89987.25
$ws.Range("H75").Value = 24500
$ws.Range("J75").Value = 24500
$ws.Range("L75").Value = 24500
$ws.Range("N75").Value = -26372
$ws.Range("H78").Value = 24500
$ws.Range("J78").Value = 24500
$ws.Range("L78").Value = 73500
$ws.Range("N78").Value = -82860
$ws.Range("H100").Value = 2547.5
$ws.Range("J100").Value = 2663.3333
$ws.Range("L100").Value = 5326.6666
$ws.Range("N100").Value = -6408.6666
$ws.Range("H132").Value = 11500.75
$ws.Range("I132").Value = 8001.5
$ws.Range("K132").Value = 24004.5
$ws.Range("M132").Value = -21474.5
